$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.357.87'
$ws.Range('E2').Value = '  -2.47%  '
$ws.Range('D3').Value = '1.775.62'
$ws.Range('E3').Value = '  -0.95%  '
$c = $ws.Range('D4')
$c.NumberFormat = '@'
$c.Value = '1.002'
$c.Style = 'Normal'
$ws.Range('E4').Value = '  -0.81%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '1.001'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  -0.71%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '305.49'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  -0.77%  '
$c = $ws.Range('D7')
$c.NumberFormat = '@'
$c.Value = '0.4231'
$c.Style = 'Normal'
$ws.Range('E7').Value = '  +1.32%  '
$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '0.3594'
$c.Style = 'Normal'
$ws.Range('E8').Value = '  +0.86%  '
$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '0.07148'
$c.Style = 'Normal'
$ws.Range('E9').Value = '  +1.43%  '
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '0.8356'
$c.Style = 'Normal'
$ws.Range('E10').Value = '  -0.88%  '
$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '20.39'
$c.Style = 'Normal'
$ws.Range('E11').Value = '  +1.78%  '
$ws.Range('D12').Value = '1.772.31'
$ws.Range('E12').Value = '  -4.55%  '
$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '6.437'
$c.Style = 'Normal'
$ws.Range('E13').Value = '  +1.54%  '
$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '5.241'
$c.Style = 'Normal'
$ws.Range('E14').Value = '  -0.26%  '
$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '0.06867'
$c.Style = 'Normal'
$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '1.005'
$c.Style = 'Normal'
$ws.Range('E16').Value = '  -0.60%  '
$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '79.02'
$c.Style = 'Normal'
$ws.Range('E17').Value = '  -0.93%  '
$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '0.000008640'
$c.Style = 'Normal'
$ws.Range('E18').Value = '  -0.85%  '
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '1.001'
$c.Style = 'Normal'
$ws.Range('E19').Value = '  -0.79%  '
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '14.88'
$c.Style = 'Normal'
$ws.Range('E20').Value = '  -0.96%  '
$ws.Range('D21').Value = '26.361.55'
$ws.Range('E21').Value = '  -3.47%  '
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '5.072'
$c.Style = 'Normal'
$ws.Range('E22').Value = '  +0.70%  '
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '10.88'
$c.Style = 'Normal'
$ws.Range('E23').Value = '  +2.19%  '
$ws.Range('D24').Value = '1.994.48'
$ws.Range('E24').Value = '  -4.02%  '
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '151.60'
$c.Style = 'Normal'
$ws.Range('E25').Value = '  -0.84%  '
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '1.788'
$c.Style = 'Normal'
$ws.Range('E26').Value = '  -9.20%  '
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '18.00'
$c.Style = 'Normal'
$ws.Range('E27').Value = '  -0.72%  '
$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '5.064'
$c.Style = 'Normal'
$ws.Range('E28').Value = '  +2.07%  '
$ws.Range('E29').Value = '  +1.58%  '
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '1.828'
$c.Style = 'Normal'
$ws.Range('E30').Value = '  +10.13%  '
$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '0.08811'
$c.Style = 'Normal'
$ws.Range('E31').Value = '  -0.83%  '
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '0.7244'
$c.Style = 'Normal'
$ws.Range('E32').Value = '  +0.10%  '
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '1.119'
$c.Style = 'Normal'
$ws.Range('E33').Value = '  +3.92%  '
$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '4.316'
$c.Style = 'Normal'
$ws.Range('E34').Value = '  -0.63%  '
$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '1.000'
$c.Style = 'Normal'
$ws.Range('E35').Value = '  -0.81%  '
$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '2.730'
$c.Style = 'Normal'
$ws.Range('E36').Value = '  -5.23%  '
$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '1.080'
$c.Style = 'Normal'
$ws.Range('E37').Value = '  +0.61%  '
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '0.05114'
$c.Style = 'Normal'
$ws.Range('E38').Value = '  +0.15%  '
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '0.01879'
$c.Style = 'Normal'
$ws.Range('E39').Value = '  -0.57%  '
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '0.4912'
$c.Style = 'Normal'
$ws.Range('E40').Value = '  -0.69%  '
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '0.1603'
$c.Style = 'Normal'
$ws.Range('E41').Value = '  -0.78%  '
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '2.603'
$c.Style = 'Normal'
$ws.Range('E42').Value = '  -2.43%  '
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '6.297'
$c.Style = 'Normal'
$ws.Range('E43').Value = '  +2.37%  '
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '7.962'
$c.Style = 'Normal'
$ws.Range('E44').Value = '  -0.45%  '
$ws.Range('B45').Value = 'Quant'
$ws.Range('C45').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '104.25'
$c.Style = 'Normal'
$ws.Range('E45').Value = '  -0.17%  '
$ws.Range('B46').Value = 'PaxDollar'
$ws.Range('C46').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '1.000'
$c.Style = 'Normal'
$ws.Range('E46').Value = '  -0.81%  '
$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '10.15'
$c.Style = 'Normal'
$ws.Range('E47').Value = '  -0.05%  '
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '1.631'
$c.Style = 'Normal'
$ws.Range('E48').Value = '  +2.80%  '
$ws.Range('E49').Value = '  -2.21%  '
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '0.4440'
$c.Style = 'Normal'
$ws.Range('E50').Value = '  -1.92%  '
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '1.718'
$c.Style = 'Normal'
$ws.Range('E51').Value = '  +3.55%  '
